$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.405.62'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  -0.08%  '
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.847.52'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  -0.20%  '
$style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.01%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.69'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -0.66%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6317'
$ws.Range('D6').Style = $style
$ws.Range('E7').Value = '  +0.01%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07572'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  +0.00%  '
$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2934'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  -0.93%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.57'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  -0.14%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07717'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -0.13%  '
$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.876.29'
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  -5.48%  '
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.003'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +0.08%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6806'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -0.60%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001046'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +5.47%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.51'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +0.59%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.131.23'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -5.88%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.175'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  -0.22%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.445.28'
$ws.Range('D19').Style = $style
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.84'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('E21').Value = '  -0.24%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = $style
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.499'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -1.48%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +0.00%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.68'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +0.43%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1395'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +0.48%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.344'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('E28').Value = '  -0.48%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.464'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  -0.38%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.300'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +3.34%  '
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05649'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -1.94%  '
$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.105'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  -0.70%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.026'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +0.11%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.853'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -0.20%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.158'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -0.08%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7071'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -1.25%  '
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.592'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  -0.13%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.250.03'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  +0.32%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.773'
$ws.Range('D40').Style = $style
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.388'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +4.85%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9034'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +0.04%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.82'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +0.12%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.94'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -1.56%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000119'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +0.78%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.104'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -0.72%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4003'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.676'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.929'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  -2.95%  '
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1124'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  -0.10%  '
